$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(26,1).Value = 25
$ws.Cells.Item(26,2).Value = 'Total Memory: 15.86 GB, Used Memory: 5.96 GB, Total Disk Space: 237.84 GB'
$ws.Cells.Item(26,3).Value = 'Sovan.Souern'
$ws.Cells.Item(26,4).Value = '1L0N1W2'
$ws.Cells.Item(26,5).Value = 'AMD64'
$ws.Cells.Item(26,6).Value = 'Windows'
$ws.Cells.Item(26,7).Value = 'PNCL114'
$ws.Cells.Item(26,8).Value = 'AT/AT COMPATIBLE'

$ws.Cells.Item(27,1).Value = 26
$ws.Cells.Item(27,2).Value = 'Total Memory: 15.86 GB, Used Memory: 6.18 GB, Total Disk Space: 237.84 GB'
$ws.Cells.Item(27,3).Value = 'Sovan.Souern'
$ws.Cells.Item(27,4).Value = '1L0N1W2'
$ws.Cells.Item(27,5).Value = 'AMD64'
$ws.Cells.Item(27,6).Value = 'Windows'
$ws.Cells.Item(27,7).Value = 'PNCL114'
$ws.Cells.Item(27,8).Value = 'AT/AT COMPATIBLE'

$ws.Cells.Item(28,1).Value = 27
$ws.Cells.Item(28,2).Value = 'Total Memory: 15.86 GB, Used Memory: 6.15 GB, Total Disk Space: 237.84 GB'
$ws.Cells.Item(28,3).Value = 'Sovan.Souern'
$ws.Cells.Item(28,4).Value = '1L0N1W2'
$ws.Cells.Item(28,5).Value = 'AMD64'
$ws.Cells.Item(28,6).Value = 'Windows'
$ws.Cells.Item(28,7).Value = 'PNCL114'
$ws.Cells.Item(28,8).Value = 'AT/AT COMPATIBLE'

$ws.Cells.Item(29,1).Value = 28
$ws.Cells.Item(29,2).Value = 'Total Memory: 15.86 GB, Used Memory: 6.10 GB, Total Disk Space: 237.84 GB'
$ws.Cells.Item(29,3).Value = 'Sovan.Souern'
$ws.Cells.Item(29,4).Value = '1L0N1W2'
$ws.Cells.Item(29,5).Value = 'AMD64'
$ws.Cells.Item(29,6).Value = 'Windows'
$ws.Cells.Item(29,7).Value = 'PNCL114'
$ws.Cells.Item(29,8).Value = 'AT/AT COMPATIBLE'

$ws.Cells.Item(30,1).Value = 29
$ws.Cells.Item(30,2).Value = 'Total Memory: 15.86 GB, Used Memory: 6.33 GB, Total Disk Space: 237.84 GB'
$ws.Cells.Item(30,3).Value = 'Sovan.Souern'
$ws.Cells.Item(30,4).Value = '1L0N1W2'
$ws.Cells.Item(30,5).Value = 'AMD64'
$ws.Cells.Item(30,6).Value = 'Windows'
$ws.Cells.Item(30,7).Value = 'PNCL114'
$ws.Cells.Item(30,8).Value = 'AT/AT COMPATIBLE'

$ws.Cells.Item(31,1).Value = 30
$ws.Cells.Item(31,2).Value = 'Total Memory: 15.86 GB, Used Memory: 5.49 GB, Total Disk Space: 237.84 GB'
$ws.Cells.Item(31,3).Value = 'Sovan.Souern'
$ws.Cells.Item(31,4).Value = '1L0N1W2'
$ws.Cells.Item(31,5).Value = 'AMD64'
$ws.Cells.Item(31,6).Value = 'Windows'
$ws.Cells.Item(31,7).Value = 'PNCL114'
$ws.Cells.Item(31,8).Value = 'AT/AT COMPATIBLE'

$ws.Cells.Item(32,1).Value = 31
$ws.Cells.Item(32,2).Value = 'Total Memory: 15.86 GB, Used Memory: 5.81 GB, Total Disk Space: 237.84 GB'
$ws.Cells.Item(32,3).Value = 'Sovan.Souern'
$ws.Cells.Item(32,4).Value = '1L0N1W2'
$ws.Cells.Item(32,5).Value = 'AMD64'
$ws.Cells.Item(32,6).Value = 'Windows'
$ws.Cells.Item(32,7).Value = 'PNCL114'
$ws.Cells.Item(32,8).Value = 'AT/AT COMPATIBLE'

$ws.Cells.Item(33,1).Value = 32
$ws.Cells.Item(33,2).Value = 'Total Memory: 15.86 GB, Used Memory: 5.64 GB, Total Disk Space: 237.84 GB'
$ws.Cells.Item(33,3).Value = 'Sovan.Souern'
$ws.Cells.Item(33,4).Value = '1L0N1W2'
$ws.Cells.Item(33,5).Value = 'AMD64'
$ws.Cells.Item(33,6).Value = 'Windows'
$ws.Cells.Item(33,7).Value = 'PNCL114'
$ws.Cells.Item(33,8).Value = 'AT/AT COMPATIBLE'

$ws.Cells.Item(34,1).Value = 33
$ws.Cells.Item(34,2).Value = 'Total Memory: 15.86 GB, Used Memory: 6.24 GB, Total Disk Space: 237.84 GB'
$ws.Cells.Item(34,3).Value = 'Sovan.Souern'
$ws.Cells.Item(34,4).Value = '1L0N1W2'
$ws.Cells.Item(34,5).Value = 'AMD64'
$ws.Cells.Item(34,6).Value = 'Windows'
$ws.Cells.Item(34,7).Value = 'PNCL114'
$ws.Cells.Item(34,8).Value = 'AT/AT COMPATIBLE'

$ws.Cells.Item(35,1).Value = 34
$ws.Cells.Item(35,2).Value = 'Total Memory: 15.86 GB, Used Memory: 6.24 GB, Total Disk Space: 237.84 GB'
$ws.Cells.Item(35,3).Value = 'Sovan.Souern'
$ws.Cells.Item(35,4).Value = '1L0N1W2'
$ws.Cells.Item(35,5).Value = 'AMD64'
$ws.Cells.Item(35,6).Value = 'Windows'
$ws.Cells.Item(35,7).Value = 'PNCL114'
$ws.Cells.Item(35,8).Value = 'AT/AT COMPATIBLE'

$ws.Cells.Item(36,1).Value = 35
$ws.Cells.Item(36,2).Value = 'Total Memory: 15.86 GB, Used Memory: 6.24 GB, Total Disk Space: 237.84 GB'
$ws.Cells.Item(36,3).Value = 'Sovan.Souern'
$ws.Cells.Item(36,4).Value = '1L0N1W2'
$ws.Cells.Item(36,5).Value = 'AMD64'
$ws.Cells.Item(36,6).Value = 'Windows'
$ws.Cells.Item(36,7).Value = 'PNCL114'
$ws.Cells.Item(36,8).Value = 'AT/AT COMPATIBLE'

$ws.Cells.Item(37,1).Value = 36
$ws.Cells.Item(37,2).Value = 'Total Memory: 15.86 GB, Used Memory: 6.26 GB, Total Disk Space: 237.84 GB'
$ws.Cells.Item(37,3).Value = 'Sovan.Souern'
$ws.Cells.Item(37,4).Value = '1L0N1W2'
$ws.Cells.Item(37,5).Value = 'AMD64'
$ws.Cells.Item(37,6).Value = 'Windows'
$ws.Cells.Item(37,7).Value = 'PNCL114'
$ws.Cells.Item(37,8).Value = 'AT/AT COMPATIBLE'

$ws.Cells.Item(38,1).Value = 37
$ws.Cells.Item(38,2).Value = 'Total Memory: 15.86 GB, Used Memory: 5.71 GB, Total Disk Space: 237.84 GB'
$ws.Cells.Item(38,3).Value = 'Sovan.Souern'
$ws.Cells.Item(38,4).Value = '1L0N1W2'
$ws.Cells.Item(38,5).Value = 'AMD64'
$ws.Cells.Item(38,6).Value = 'Windows'
$ws.Cells.Item(38,7).Value = 'PNCL114'
$ws.Cells.Item(38,8).Value = 'AT/AT COMPATIBLE'

$ws.Cells.Item(39,1).Value = 38
$ws.Cells.Item(39,2).Value = 'Total Memory: 15.86 GB, Used Memory: 7.50 GB, Total Disk Space: 237.84 GB'
$ws.Cells.Item(39,3).Value = 'Sovan.Souern'
$ws.Cells.Item(39,4).Value = '1L0N1W2'
$ws.Cells.Item(39,5).Value = 'AMD64'
$ws.Cells.Item(39,6).Value = 'Windows'
$ws.Cells.Item(39,7).Value = 'PNCL114'
$ws.Cells.Item(39,8).Value = 'AT/AT COMPATIBLE'

$ws.Cells.Item(40,1).Value = 39
$ws.Cells.Item(40,2).Value = 'Total Memory: 15.86 GB, Used Memory: 7.55 GB, Total Disk Space: 237.84 GB'
$ws.Cells.Item(40,3).Value = 'Sovan.Souern'
$ws.Cells.Item(40,4).Value = '1L0N1W2'
$ws.Cells.Item(40,5).Value = 'AMD64'
$ws.Cells.Item(40,6).Value = 'Windows'
$ws.Cells.Item(40,7).Value = 'PNCL114'
$ws.Cells.Item(40,8).Value = 'AT/AT COMPATIBLE'

$ws.Cells.Item(41,1).Value = 40
$ws.Cells.Item(41,2).Value = 'Total Memory: 15.86 GB, Used Memory: 5.84 GB, Total Disk Space: 237.84 GB'
$ws.Cells.Item(41,3).Value = 'Sovan.Souern'
$ws.Cells.Item(41,4).Value = '1L0N1W2'
$ws.Cells.Item(41,5).Value = 'AMD64'
$ws.Cells.Item(41,6).Value = 'Windows'
$ws.Cells.Item(41,7).Value = 'PNCL114'
$ws.Cells.Item(41,8).Value = 'AT/AT COMPATIBLE'

# Copy formatting (style s="2") from the last pre-existing data row (25) down to the new rows
$ws.Range("A25:H25").Copy()
$ws.Range("A26:H41").PasteSpecial(-4122)
$excel.CutCopyMode = 0
